$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Oct 30 2020", " Abu Dhabi", "Royals won by 7 wickets (with 15 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Mandeep Singh ", "'0", "'1", "'0", "'0", "'0.00"),
    @(" Oct 4 2020", " Dubai (DSC)", "Super Kings won by 10 wickets (with 14 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "Mandeep Singh ", "'27", "'16", "'0", "'2", "'168.75"),
    @(" Oct 24 2020", " Dubai (DSC)", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Mandeep Singh ", "'17", "'14", "'1", "'0", "'121.42"),
    @(" Oct 26 2020", " Sharjah", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kings XI Punjab", "Kolkata Knight Riders", "Mandeep Singh ", "'66", "'56", "'8", "'2", "'117.85"),
    @(" Oct 10 2020", " Abu Dhabi", "KKR won by 2 runs", "Kings XI Punjab", "Kolkata Knight Riders", "Mandeep Singh ", "'0", "'1", "'0", "'0", "'0.00"),
    @(" Oct 8 2020", " Dubai (DSC)", "Sunrisers won by 69 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Mandeep Singh ", "'6", "'6", "'0", "'0", "'100.00")
)

$rowIndex = 3
foreach ($rowData in $data) {
    $colIndex = 1
    foreach ($val in $rowData) {
        $cell = $ws.Cells.Item($rowIndex, $colIndex)
        $cell.Value = $val
        $cell.Style = "Normal"
        $colIndex++
    }
    $rowIndex++
}
